# The workbook's "vocabulary" sheet lists PREFIX rows starting at row 3.
# A new PREFIX row for "xsd" (http://www.w3.org/2001/XMLSchema#) needs to be
# inserted right before the existing "qudt" row (currently row 10), pushing
# qudt/unit and everything below it down by one row. The "dct:modified"
# timestamp (now shifted to row 21) also gets bumped to reflect the
# regeneration time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10 ("qudt" PREFIX row); this shifts
# rows 10:23 down to 11:24 and extends the sheet dimension automatically.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the "xsd" prefix entry.
$ws.Range("A10").Value = "PREFIX"
$ws.Range("B10").Value = "xsd"
$ws.Range("C10").Value = "http://www.w3.org/2001/XMLSchema#"

# Match the sheet's convention (every other data row spans A:AM with empty
# string placeholders in the unused trailing columns) by touching D10:AM10.
$ws.Range("D10:AM10").Value = ""

# Update the "dct:modified" timestamp (shifted down to row 21) to the new
# generation time recorded for this export.
$ws.Range("B21").Value = "2023-01-29T20:06:06+00:00"
